# Generate Report for Handback
# Update status for e2fbda3d-... row from "Ready for handoff" to
# "Handed back: in sync with en-US" on all sheets, and refresh the
# "Latest Handback DateTime" values on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $statusText
$zhcn.Range("G2").Value = "2016-03-04 11:26:15"
$zhcn.Range("G3").Value = "2016-03-04 11:26:15"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $statusText
$dede.Range("G2").Value = "2016-03-04 11:26:41"
$dede.Range("G3").Value = "2016-03-04 11:26:41"
